$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.855.81"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "3.537.33"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.34%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  +4.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.640"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000308"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "4.102.50"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "69.959.04"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "584.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.65"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.85"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.47%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.528.72"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.991"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "63.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +17.99%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "531.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.08%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.93%  "
$ws.Range("D40").Value = "0.0₃0777"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "3.530.84"
$ws.Range("E41").Value = "  +5.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.135"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "
$ws.Range("E44").Value = "  +1.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("E46").Value = "  -5.48%  "
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  -3.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.01%  "
